$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 336.32144
$ws.Range("I28").Value = 176.29167
$ws.Range("J28").Value = 1296.5
$ws.Range("K28").Value = 176.29167
$ws.Range("L28").Value = 1296.5
$ws.Range("M28").Value = 308.70833
$ws.Range("N28").Value = -2266.5

$ws.Range("H107").Value = 48247.617
$ws.Range("I107").Value = 71968.57000000001
$ws.Range("J107").Value = 805.7143
$ws.Range("K107").Value = 71968.57000000001
$ws.Range("L107").Value = 805.7143
$ws.Range("M107").Value = -70048.57000000001
$ws.Range("N107").Value = -4645.7143

$ws.Range("H132").Value = 2802.1135
$ws.Range("I132").Value = 1918.7142
$ws.Range("J132").Value = 4348.0625
$ws.Range("K132").Value = 5756.142599999999
$ws.Range("L132").Value = 13044.1875
$ws.Range("M132").Value = -3226.142599999999
$ws.Range("N132").Value = -18104.1875

$ws.Range("H138").Value = 36038676
$ws.Range("I138").Value = 43480010
$ws.Range("J138").Value = 23813630
$ws.Range("K138").Value = 130440030
$ws.Range("L138").Value = 71440890
$ws.Range("M138").Value = -130434890
$ws.Range("N138").Value = -71451170

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 41534.145
$ws.Range("J130").Value = 41534.145
$ws.Range("L130").Value = 41534.145
$ws.Range("N130").Value = -51574.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1379.4615
$ws.Range("I16").Value = 1243.3
$ws.Range("J16").Value = 1833.3334
$ws.Range("K16").Value = 1243.3
$ws.Range("L16").Value = 1833.3334
$ws.Range("M16").Value = -956.3
$ws.Range("N16").Value = -2407.3334

$ws.Range("H31").Value = 3290
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3290
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3290
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -3880

$ws.Range("H34").Value = 3290
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3290
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3290
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3694

$ws.Range("H113").Value = 1379.4615
$ws.Range("I113").Value = 1243.3
$ws.Range("J113").Value = 1833.3334
$ws.Range("K113").Value = 1243.3
$ws.Range("L113").Value = 1833.3334
$ws.Range("M113").Value = 926.7
$ws.Range("N113").Value = -6173.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 664.6875
$ws.Range("I5").Value = 409.33334
$ws.Range("J5").Value = 1152.1818
$ws.Range("K5").Value = 1228.00002
$ws.Range("L5").Value = 3456.5454
$ws.Range("M5").Value = -1116.00002
$ws.Range("N5").Value = -3680.5454

$ws.Range("H22").Value = 1928.5714
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1928.5714
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 5785.7142
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -6123.7142

$ws.Range("H27").Value = 1928.5714
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1928.5714
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 5785.7142
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -5989.7142

$ws.Range("H113").Value = 504.1887
$ws.Range("I113").Value = 444.19354
$ws.Range("J113").Value = 588.7273
$ws.Range("K113").Value = 1332.58062
$ws.Range("L113").Value = 1766.1819
$ws.Range("M113").Value = 837.41938
$ws.Range("N113").Value = -6106.1819

$ws.Range("H127").Value = 568
$ws.Range("J127").Value = 568
$ws.Range("L127").Value = 1704
$ws.Range("N127").Value = -11624

$ws.Range("H135").Value = 664.6875
$ws.Range("I135").Value = 409.33334
$ws.Range("J135").Value = 1152.1818
$ws.Range("K135").Value = 3684.00006
$ws.Range("L135").Value = 10369.6362
$ws.Range("M135").Value = -1149.00006
$ws.Range("N135").Value = -15439.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 4211.143
$ws.Range("I126").Value = 2199.4736
$ws.Range("J126").Value = 6600
$ws.Range("K126").Value = 6598.4208
$ws.Range("L126").Value = 19800
$ws.Range("M126").Value = -4128.4208
$ws.Range("N126").Value = -24740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 653.8570999999999
$ws.Range("I55").Value = 435.57144
$ws.Range("J55").Value = 763
$ws.Range("K55").Value = 435.57144
$ws.Range("L55").Value = 763
$ws.Range("M55").Value = -262.57144
$ws.Range("N55").Value = -1109

$ws.Range("H61").Value = 1364.3125
$ws.Range("I61").Value = 1361.4375
$ws.Range("J61").Value = 1367.1875
$ws.Range("K61").Value = 1361.4375
$ws.Range("L61").Value = 1367.1875
$ws.Range("M61").Value = -1159.4375
$ws.Range("N61").Value = -1771.1875

$ws.Range("H113").Value = 1364.3125
$ws.Range("I113").Value = 1361.4375
$ws.Range("J113").Value = 1367.1875
$ws.Range("K113").Value = 1361.4375
$ws.Range("L113").Value = 1367.1875
$ws.Range("M113").Value = 808.5625
$ws.Range("N113").Value = -5707.1875

$ws.Range("H122").Value = 3118.0356
$ws.Range("I122").Value = 2573.611
$ws.Range("J122").Value = 4098
$ws.Range("K122").Value = 7720.833
$ws.Range("L122").Value = 12294
$ws.Range("M122").Value = -5270.833
$ws.Range("N122").Value = -17194

$ws.Range("H132").Value = 2286.9355
$ws.Range("I132").Value = 2012.4783
$ws.Range("J132").Value = 3076
$ws.Range("K132").Value = 6037.4349
$ws.Range("L132").Value = 9228
$ws.Range("M132").Value = -3507.4349
$ws.Range("N132").Value = -14288

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 833923
$ws.Range("I100").Value = 583.4286
$ws.Range("J100").Value = 2000598.4
$ws.Range("K100").Value = 1166.8572
$ws.Range("L100").Value = 4001196.8
$ws.Range("M100").Value = -625.8571999999999
$ws.Range("N100").Value = -4002278.8

$ws.Range("H132").Value = 4057.2415
$ws.Range("I132").Value = 4109.1
$ws.Range("J132").Value = 3942
$ws.Range("K132").Value = 12327.3
$ws.Range("L132").Value = 11826
$ws.Range("M132").Value = -9797.300000000001
$ws.Range("N132").Value = -16886

Write-Output "Applied all changes"